# excel_processor.py able to process multiple sheets and giving multiple
# excel as output: insert a fresh blank "Sheet1" ahead of the existing
# "High Impact" data sheet, tidy a couple of data values on the data sheet,
# append a free-form description row, and grow the table to match.

$wb = $excel.ActiveWorkbook

# The workbook currently has a single sheet ("High Impact"). Insert a new,
# blank worksheet in front of it - Worksheets.Add() with no args inserts
# before the active sheet and Excel names it "Sheet1".
$newSheet = $wb.Worksheets.Add()

# Grab a handle to the original data sheet by name (now pushed to position 2).
$dataSheet = $wb.Worksheets.Item("High Impact")

# K2/L2 held verbose "YYYY-MM-DD 00:00:00" text; trim to a plain date string.
# Force the cell to Text first so Excel doesn't reinterpret the literal as a
# real date serial, then drop the number-format override so no stray style
# is left behind on the cell.
$dataSheet.Range("K2").NumberFormat = "@"
$dataSheet.Range("K2").Value = "2025-09-26"
$dataSheet.Range("K2").ClearFormats()

$dataSheet.Range("L2").NumberFormat = "@"
$dataSheet.Range("L2").Value = "2025-09-26"
$dataSheet.Range("L2").ClearFormats()

# U3/V3 were live SUM() formulas referencing the single data row; replace
# them with their resolved literal values (same trick to keep them as text,
# matching the source values stored in U2/V2).
$dataSheet.Range("U3").NumberFormat = "@"
$dataSheet.Range("U3").Value = "13813169"
$dataSheet.Range("U3").ClearFormats()

$dataSheet.Range("V3").NumberFormat = "@"
$dataSheet.Range("V3").Value = "203744"
$dataSheet.Range("V3").ClearFormats()

# Append the free-form description note in a new row 4.
$dataSheet.Range("A4").Value = "Free-Form Description: Please use this space to fill in details specific to target audience and any special needs for campaign execution such as minimum lead time for creative deliverables, other available takeover opportunities or potential beta test/first to market details"

# Grow the table (and its autofilter) so row 4 is included in its range.
$lo = $dataSheet.ListObjects.Item(1)
$lo.Resize($dataSheet.Range("A1:Y4"))
